$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell values taken from the authoritative diff (cryptos.xlsx update).
$updates = [ordered]@{
    "D2" = "64.809.78"
    "E2" = "  +1.61%  "
    "D3" = "3.154.60"
    "E3" = "  +0.57%  "
    "D4" = "1.01"
    "E4" = "  +0.65%  "
    "D5" = "592.46"
    "E5" = "  +0.98%  "
    "D6" = "153.15"
    "E6" = "  +4.70%  "
    "E7" = "  +0.32%  "
    "D8" = "3.153.45"
    "E8" = "  +0.61%  "
    "D9" = "0.536"
    "E9" = "  +1.29%  "
    "D10" = "0.162"
    "E10" = "  +0.33%  "
    "D11" = "5.99"
    "E11" = "  +4.12%  "
    "D12" = "0.466"
    "E12" = "  +1.77%  "
    "D13" = "38.67"
    "E13" = "  +4.93%  "
    "D14" = "0.0000248"
    "E14" = "  +0.56%  "
    "D15" = "3.680.88"
    "E15" = "  +0.79%  "
    "E16" = "  -0.27%  "
    "D17" = "7.32"
    "E17" = "  +3.38%  "
    "D18" = "64.459.25"
    "E18" = "  +1.38%  "
    "D19" = "3.157.73"
    "E19" = "  +0.80%  "
    "D20" = "474.74"
    "E20" = "  +2.24%  "
    "D21" = "15.01"
    "E21" = "  +5.26%  "
    "D22" = "0.751"
    "E22" = "  +2.35%  "
    "D23" = "7.70"
    "E23" = "  +3.74%  "
    "D24" = "13.52"
    "E24" = "  +4.60%  "
    "D25" = "2.41"
    "E25" = "  +9.66%  "
    "D26" = "82.36"
    "E26" = "  +1.52%  "
    "B27" = "Dai"
    "C27" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D27" = "1.00"
    "E27" = "  -0.04%  "
    "B28" = "RenderToken"
    "C28" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D28" = "9.97"
    "E28" = "  +7.46%  "
    "D29" = "2.74"
    "E29" = "  +2.45%  "
    "D30" = "7.43"
    "E30" = "  +6.27%  "
    "D31" = "2.23"
    "E31" = "  +0.91%  "
    "D32" = "1.01"
    "E32" = "  +0.60%  "
    "D33" = "0.118"
    "E33" = "  +7.09%  "
    "D34" = "27.81"
    "E34" = "  +3.13%  "
    "D35" = "0.0₃0878"
    "E35" = "  +3.97%  "
    "D36" = "3.54"
    "E36" = "  +6.52%  "
    "E37" = "  +2.89%  "
    "D38" = "6.22"
    "E38" = "  +3.52%  "
    "D39" = "2.32"
    "E39" = "  +1.18%  "
    "D40" = "466.50"
    "E40" = "  +5.92%  "
    "D41" = "9.37"
    "E41" = "  +6.42%  "
    "D42" = "51.39"
    "E42" = "  +0.23%  "
    "E43" = "  +8.69%  "
    "D44" = "0.0382"
    "E44" = "  +3.16%  "
    "D45" = "2.900.05"
    "E45" = "  -0.61%  "
    "E46" = "  +3.50%  "
    "D47" = "38.82"
    "E47" = "  +5.11%  "
    "D48" = "131.43"
    "E48" = "  +3.68%  "
    "D49" = "25.97"
    "E49" = "  +7.62%  "
    "D50" = "2.30"
    "E50" = "  +5.89%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. "1.01", "592.46")
    # are not silently coerced to doubles / lose trailing zeros.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop the explicit text format again so the cell matches the sheet's
    # default (unstyled) data cells, as in the original workbook.
    $cell.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
